$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45436

$ws.Range("D14").Value = 83.175
$ws.Range("D15").Value = 108.235
$ws.Range("D16").Value = 129.622
$ws.Range("D17").Value = 207.394
$ws.Range("D18").Value = 259.243
$ws.Range("D19").Value = 324.054
$ws.Range("D20").Value = 360.78
$ws.Range("D21").Value = 399.666
